$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.156126260757446
$ws.Range("B1").Value = 2.533568382263184
$ws.Range("C1").Value = 2.675328969955444
$ws.Range("D1").Value = 3.264136075973511
$ws.Range("E1").Value = 2.319345235824585
